$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A:B").Insert()

$ws.Range("A1").Value = "Country"
$ws.Range("B1").Value = "Chronosequence"

$ws.Range("A2:A6").Value = "CRI"
$ws.Range("B2:B6").Value = "test"

$ws.Columns.Item(2).ColumnWidth = 16.140625

[void]$ws.Range("B2:B6").Select()

Write-Output "done"
